# Actualización automática 2025-06-02 14:06:09
# Add a new "PRESUPUESTO" (budget) column G to the "VENTA MENSUAL" sheet,
# mirroring the formatting of the existing "junio" column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column width. Excel's ColumnWidth (character units) reads back ~0.8333
# less than the width persisted in the XML, so offset to land on 17 exactly.
$ws.Columns.Item(7).ColumnWidth = 17 - 5/6

# Header cell G1 ("PRESUPUESTO"), formatted like the other header cells (F1).
$ws.Cells.Item(1, 6).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# Data cells G2:G4 (value 0), formatted like the corresponding row in column F.
$ws.Cells.Item(2, 6).Copy()
$ws.Cells.Item(2, 7).PasteSpecial(-4122)
$ws.Cells.Item(2, 7).Value = 0

$ws.Cells.Item(3, 6).Copy()
$ws.Cells.Item(3, 7).PasteSpecial(-4122)
$ws.Cells.Item(3, 7).Value = 0

$ws.Cells.Item(4, 6).Copy()
$ws.Cells.Item(4, 7).PasteSpecial(-4122)
$ws.Cells.Item(4, 7).Value = 0

$excel.CutCopyMode = $false
